$d = $word.ActiveDocument

# Locate the "Μέθοδοι Εργασίας" Heading-1 paragraph. The paragraph we need
# to edit is the empty, directly-formatted (bold / Century Gothic /
# "heading 1" pStyle) paragraph sitting 3 paragraphs before it. Anchoring
# on this stable text is more robust than a hard-coded paragraph index.
$count = $d.Paragraphs.Count
$anchorIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -eq ("Μέθοδοι Εργασίας" + [char]13)) {
        $anchorIndex = $i
    }
}

$target = $d.Paragraphs.Item($anchorIndex - 3)

# Sanity-check we found the right, currently-empty "Heading 1" paragraph
# before touching it.
if ($target.Style.NameLocal -eq "Heading 1" -and $target.Range.Text -eq [char]13) {
    # Strip it back down to a completely bare/default paragraph (same as
    # its now-plain neighbours): re-styling it to Normal clears the
    # explicit pStyle, and resetting the font clears the direct
    # paragraph-mark character formatting (bold, Century Gothic font,
    # explicit color) that rode along with it.
    $target.Style = $d.Styles.Item("Normal")
    $target.Range.Font.Reset() | Out-Null
}
